# Generate Report for Handback
# Marks the localization status report as "handed back" for each locale,
# filling in the Latest Target File / Latest Handback File / Latest
# Handback DateTime columns and updating the Status text everywhere it
# appears (Overview summary sheet + each locale sheet).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdName = "b0dcb24c-3f80-4bba-a210-45e5c733804f.md"

# --- Overview sheet: just the rolled-up status text -----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- Per-locale sheets -----------------------------------------------------
# locale sheet name -> [xlf file name, handback datetime text]
$locales = @(
    @{ Name = "zh-cn"; Xlf = "b0dcb24c-3f80-4bba-a210-45e5c733804f.baa21ee89ecc3fefd816f6131e1f61b1bd080ca5.zh-cn.xlf"; HandbackTime = "2016-03-19 10:49:37" },
    @{ Name = "de-de"; Xlf = "b0dcb24c-3f80-4bba-a210-45e5c733804f.baa21ee89ecc3fefd816f6131e1f61b1bd080ca5.de-de.xlf"; HandbackTime = "2016-03-19 10:49:42" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)
    $xlfName = $locale.Xlf
    $handbackTime = $locale.HandbackTime

    # Snapshot the existing hyperlink addresses before we rebuild the
    # collection, so the new Latest Target File / Latest Handback File
    # links point at the same external targets as the existing md / xlf
    # hyperlinks on each row.
    $addr = @{}
    $disp = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $a1 = $hl.Range.Address()
        $addr[$a1] = $hl.Address
        $disp[$a1] = $hl.TextToDisplay
    }

    $mdAddr2 = $addr['$A$2']
    $xlfAddr2 = $addr['$D$2']
    $mdAddr3 = $addr['$A$3']
    $xlfAddr3 = $addr['$D$3']

    # Status column
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # New columns: Latest Target File (F), Latest Handback File (G)
    $ws.Range("F2").Value = $mdName
    $ws.Range("G2").Value = $xlfName
    $ws.Range("F3").Value = $mdName
    $ws.Range("G3").Value = $xlfName

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $handbackTime
    $ws.Range("H3").Value = $handbackTime

    # Rebuild every hyperlink on the sheet so the new F2/G2/F3/G3 links
    # land in row order right after the existing row's links, matching
    # how the handback job emits them (row by row).
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdAddr2, "", "", $disp['$A$2']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $mdAddr2, "", "", $disp['$B$2']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $xlfAddr2, "", "", $disp['$D$2']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdAddr2, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfAddr2, "", "", $xlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdAddr3, "", "", $disp['$A$3']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdAddr3, "", "", $disp['$B$3']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $xlfAddr3, "", "", $disp['$D$3']) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdAddr3, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfAddr3, "", "", $xlfName) | Out-Null
}
